$word.UserName = "Sophie Breitbart"
$word.UserInitials = "SB"
$d = $word.ActiveDocument

$r0 = $d.Range(53, 86)
$c0 = $d.Comments.Add($r0, "Change to explain PVE for starred variables from lmers")

$r1 = $d.Range(53, 86)
$c1 = $d.Comments.Add($r1, "And change PVEs")
